$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mentors")
$ws.Activate()
$ws.Range("E2:E12").Value = "YES"
$ws.Rows("2:12").RowHeight = 15
